{"js": "const body = context.document.body;\n\n// Locate the paragraph that contains the sentence we need to update. The\n// original text is a single run inside a single paragraph, so a body-wide\n// search is sufficient and unambiguous.\nconst searchText =\n  \"Labeled Property Graph annotations example. Augments Reference Model. Statement Context aggregate SPO annotations:\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to update.\");\n}\n\nconst newText =\n  \"Labeled Property Graph annotations example. Augments Reference Model. Statement example, Statement context aggregates SPO annotations (Statement occurrences data in other Reference Model layers contexts):\";\n\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Labeled Property Graph annotations example. Augments Reference Model. Statement Context aggregate SPO annotations:\"\n$newText = \"Labeled Property Graph annotations example. Augments Reference Model. Statement example, Statement context aggregates SPO annotations (Statement occurrences data in other Reference Model layers contexts):\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1\n$find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n"}
